# Aula 37 - Listando departamentos
# Adds two new note rows (31 and 32) for lesson 37 under the existing
# session "8. Departamento: Controller & View".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Text content (kept in single-quoted / here-strings so that PowerShell
# does not try to interpolate the "$" characters that appear in the
# Thymeleaf snippets).
# ---------------------------------------------------------------------

$sessionName = '8. Departamento: Controller & View'
$lessonName  = '37. Listando departamentos'

$note31 = @'
1:28
aplicando "for each" no thymeleaf para preencher tabelas HTML:
th:each="nomeVariavelQlqr : ${nomeVariavelIdenticaAoController}"
nomeVariavelQlqr: representa a linha referente a posição do for que está
nomeVariavelIdenticaAoController: será a variável que conterá a lista recebida do controller/backend
'@

$note32Run1 = @'
4:03
abordado um dos três conceitos principais no 
'@

$note32Run2 = 'Spring MVC - o Model  (os 3 são Model , ModelMap e o ModelAndView)'

$note32Run3 = ' para renderizar e construir páginas HTML com dados do controller.'

$note32 = $note32Run1 + $note32Run2 + $note32Run3

# ---------------------------------------------------------------------
# Row 31 : aula 37 / sessão 8 / nome da aula / observação (for-each note)
# ---------------------------------------------------------------------

$ws.Range("B31").Value = 37
$ws.Range("B31").Font.Color = 0

$ws.Range("C31").Value = $sessionName
$ws.Range("C31").Font.Color = 0

$ws.Range("D31").Value = $lessonName
$ws.Range("D31").Font.Color = 0
$ws.Range("D31").WrapText = $true

$ws.Range("E31").Value = $note31
$ws.Range("E31").WrapText = $true

$ws.Rows(31).RowHeight = 120

# ---------------------------------------------------------------------
# Row 32 : aula 37 / sessão 8 / nome da aula / observação (Model note)
# ---------------------------------------------------------------------

$ws.Range("B32").Value = 37
$ws.Range("B32").Font.Color = 0

$ws.Range("C32").Value = $sessionName
$ws.Range("C32").Font.Color = 0

$ws.Range("D32").Value = $lessonName
$ws.Range("D32").Font.Color = 0
$ws.Range("D32").WrapText = $true

$ws.Range("E32").Value = $note32
$ws.Range("E32").WrapText = $true

# apply rich-text (bold + red) formatting to the "Spring MVC - o Model..."
# portion of the note, matching the highlighted segment from the source.
$run2Start = $note32Run1.Length + 1
$run2Len   = $note32Run2.Length
$chars2 = $ws.Range("E32").Characters($run2Start, $run2Len)
$chars2.Font.Bold = $true
$chars2.Font.Color = 255

$run3Start = $run2Start + $run2Len
$run3Len   = $note32Run3.Length
$chars3 = $ws.Range("E32").Characters($run3Start, $run3Len)
$chars3.Font.Bold = $false
$chars3.Font.Color = 0

$ws.Rows(32).RowHeight = 60

# ---------------------------------------------------------------------
# Update the view state so that the active selection matches what was
# recorded after the edit.
# ---------------------------------------------------------------------

$ws.Range("D37").Select()
